$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 911
$ws.Range("I2").Value = 2531
$ws.Range("J2").Value = 10109
$ws.Range("K2").Value = 54
$ws.Range("L2").Value = 2731
$ws.Range("M2").Value = 158
$ws.Range("N2").Value = 1831
$ws.Range("O2").Value = 6
$ws.Range("P2").Value = 42
$ws.Range("Q2").Value = 20
$ws.Range("R2").Value = 130
$ws.Range("S2").Value = 1095
$ws.Range("T2").Value = 1731
$ws.Range("U2").Value = 148
$ws.Range("V2").Value = 15644
$ws.Range("W2").Value = 7
$ws.Range("X2").Value = 15714
$ws.Range("Y2").Value = 23
$ws.Range("Z2").Value = 215
$ws.Range("AA2").Value = 92
